$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

# New bug-log entry in row 3 (45028 = 2023-04-12, same style as A2)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 45028

$ws.Range("B3").Value = "front end navigation"
$ws.Range("C3").Value = "Front end routing was not working properly. "
$ws.Range("D3").Value = "Solved"
$ws.Range("E3").Value = "Created a custom history component that can be used outside components and hooks."
$ws.Range("F3").Value = "export const history = {`n    navigate: null,`n    location: null`n}; then initialize this in app like, history.navigate = useNavigate()"

# Row height for the wrapped content
$ws.Range("A3:F3").RowHeight = 72

# Column B needs to widen to fit the new text, matching Excel's autofit behaviour
$ws.Columns.Item(2).ColumnWidth = 16.83

# Move the active selection like the saved file shows
$ws.Range("G3").Select()
